$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trip 1 - Parkes preharvest 2020
$ws.Range("G2:G7").Value = "preharvest"

# Trip 2 - Parkes postharvest 2020
$ws.Range("G8:G13").Value = "postharvest"

# Trip 3 - Parkes PreRoll1 (pre-rolling Trip 1)
$ws.Range("G14:G19").Value = "pre-rolling"

# Trip 4 - Parkes PostRoll1 (post-rolling)
$ws.Range("G20:G22").Value = "post-rolling"

# Trip 4 - Parkes PreRoll2 2021 (pre-rolling)
$ws.Range("G23:G25").Value = "pre-rolling"

# Trip 5 - Parkes PostRoll2 2021 (post-rolling)
$ws.Range("G26:G31").Value = "post-rolling"

# New "Manipulation" column (G) header - added last so it lands after the
# treatment labels in the shared-strings table, matching the saved order
$ws.Range("G1").Value = "Manipulation"

# Match the new column width used for the Manipulation column
$ws.Columns.Item(7).ColumnWidth = 19

# Update the active selection shown in the sheet view
$ws.Range("J7").Select()
